# Applies the "South Korea K3 League" update: a handful of match rows on
# the same date got re-sorted (their data swapped/rotated between row
# positions), while row numbers / the leading index column (A) stay put.
#
# Each group below lists row numbers that form a rotation cycle: the data
# that currently lives in row group[i] ends up in row group[i+1] (wrapping
# around). Column A (the running match index, 0,1,2,...) is intentionally
# left untouched - only columns B:AD (id .. PL_AhUnder) move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycles = @(
    @(17, 18),
    @(49, 50),
    @(56, 57),
    @(93, 94),
    @(103, 104),
    @(110, 111),
    @(118, 119),
    @(121, 123, 122),
    @(124, 126, 127),
    @(129, 130, 131, 132),
    @(142, 143),
    @(147, 148),
    @(156, 158),
    @(171, 172),
    @(175, 176, 177),
    @(199, 201),
    @(207, 208)
)

foreach ($cycle in $cycles) {
    # Snapshot the B:AD values for every row in this cycle before writing
    # anything back (the ranges overlap the writes, so we must capture
    # everything up front).
    $snapshots = @()
    foreach ($r in $cycle) {
        $rng = $ws.Range("B$r`:AD$r")
        $snapshots += , $rng.Value2
    }

    $n = $cycle.Length
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        # Row at position i receives what used to be at position i+1
        # (wrapping around) - i.e. a left-rotation of the snapshots.
        $srcSnapshot = $snapshots[($i + 1) % $n]
        $ws.Range("B$destRow`:AD$destRow").Value2 = $srcSnapshot
    }
}
